$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 ---
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "2026-02-12T10:27:25.465498+00:00"
$ws.Range("E2").Value = "Sure thanks"
$ws.Range("G2").Value = "No worries bro"

# --- Update row 3 ---
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "2026-02-12T10:27:27.290161+00:00"
$ws.Range("E3").Value = "Ok thank you"
$ws.Range("G3").Value = "No worries"

# --- Add row 5 ---
$ws.Range("A5").Value = 5163876201
$ws.Range("B5").Value = "VaLaK_DEMON"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "2026-02-12T10:27:08.760567+00:00"
$ws.Range("H5").Value = $false
$ws.Range("I5").Value = $false
$ws.Range("K5").Value = "Added during extraction"

# --- Add row 6 ---
$ws.Range("A6").Value = 1900918712
$ws.Range("B6").Value = "lunaticbeast12"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = "2026-02-12T10:27:10.078943+00:00"
$ws.Range("H6").Value = $false
$ws.Range("I6").Value = $false
$ws.Range("K6").Value = "Added during extraction"
